$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'
$w14Ns = 'http://schemas.microsoft.com/office/word/2010/wordml'

# ---------------------------------------------------------------------------
# 1) Expand the lone "Df.head()" paragraph (right after "LIVE DEMO") into the
#    full "creating a DF" walkthrough block, ending with the original
#    "Df.head()" text (now proofed/split into runs).
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$null = $rng1.Find.Execute("Df.head()", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headPara = $rng1.Paragraphs(1)
$headRange = $headPara.Range

$newBlock = @"
<w:p xmlns:w="$wNs" xmlns:w14="$w14Ns" w14:paraId="2B09E7D9" w14:textId="481BBFD2" w:rsidR="00911D23" w:rsidRDefault="00911D23"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Let's start by creating a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>DF</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p xmlns:w="$wNs" xmlns:w14="$w14Ns"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&gt;&gt;&gt; d = {'col1': [1, 2], 'col2': [3, 4]}</w:t></w:r></w:p><w:p xmlns:w="$wNs" xmlns:w14="$w14Ns"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">&gt;&gt;&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>df</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>pd.DataFrame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(data=d)</w:t></w:r></w:p><w:p xmlns:w="$wNs" xmlns:w14="$w14Ns"><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">&gt;&gt;&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>df</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="$wNs" xmlns:w14="$w14Ns"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p xmlns:w="$wNs" xmlns:w14="$w14Ns"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p xmlns:w="$wNs" xmlns:w14="$w14Ns"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Df.head</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>()</w:t></w:r></w:p>
"@

$headRange.InsertXML($newBlock)

# ---------------------------------------------------------------------------
# 2) Add a <w:lastRenderedPageBreak/> in front of the "SELECT A COLUMN" run.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$null = $rng2.Find.Execute("SELECT A COLUMN", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$selPara = $rng2.Paragraphs(1)
$selRange = $selPara.Range
$selXml = @"
<w:p xmlns:w="$wNs" xmlns:w14="$w14Ns" w14:paraId="76344C1A" w14:textId="652EC782" w:rsidR="00C53CDC" w:rsidRDefault="00C53CDC" w:rsidP="00911D23"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>SELECT A COLUMN</w:t></w:r></w:p>
"@
$selRange.InsertXML($selXml)

# ---------------------------------------------------------------------------
# 3) Remove the <w:lastRenderedPageBreak/> that currently sits in front of
#    the "raw_df" run inside the "TRANSFORMING DATA" code block. The OM's
#    WordOpenXML getter already omits lastRenderedPageBreak, so re-inserting
#    a paragraph's own OOXML (round-tripped through the getter) strips it.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$null = $rng3.Find.Execute("raw_df['Winner'] = np", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$transPara = $rng3.Paragraphs(1)
$transRange = $transPara.Range
$fullPkg = $transRange.WordOpenXML
$match = [regex]::Match($fullPkg, '(?s)<w:p\b[^>]*>.*</w:p>')
$fragment = $match.Value
$fragment = $fragment -replace '^<w:p ', ('<w:p xmlns:w="' + $wNs + '" xmlns:w14="' + $w14Ns + '" ')
$transRange.InsertXML($fragment)
